$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.778.10"
$ws.Range("E2").Value = "  -5.20%  "
$ws.Range("D3").Value = "3.213.77"
$ws.Range("E3").Value = "  -8.74%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.84"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.86"
$ws.Range("E6").Value = "  -12.68%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.203.47"
$ws.Range("E8").Value = "  -8.91%  "
$ws.Range("E9").Value = "  -11.36%  "
$ws.Range("E10").Value = "  -12.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  -10.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  -16.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.99"
$ws.Range("E13").Value = "  -15.88%  "
$ws.Range("E14").Value = "  -12.38%  "
$ws.Range("D15").Value = "3.736.32"
$ws.Range("E15").Value = "  -8.75%  "
$ws.Range("D16").Value = "66.785.54"
$ws.Range("E16").Value = "  -5.20%  "
$ws.Range("D17").Value = "3.217.03"
$ws.Range("E17").Value = "  -8.66%  "
$ws.Range("E18").Value = "  -4.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "530.11"
$ws.Range("E19").Value = "  -13.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.13"
$ws.Range("E20").Value = "  -14.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  -15.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.759"
$ws.Range("E22").Value = "  -14.07%  "
$ws.Range("E23").Value = "  -12.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.93"
$ws.Range("E24").Value = "  -11.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.70"
$ws.Range("E25").Value = "  -13.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.18"
$ws.Range("E27").Value = "  -14.63%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  -15.39%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.14"
$ws.Range("E29").Value = "  -10.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.11"
$ws.Range("E30").Value = "  -15.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  -9.71%  "
$ws.Range("E32").Value = "  -9.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "544.90"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.47"
$ws.Range("E34").Value = "  -19.95%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("E35").Value = "  -16.41%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  -5.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0426"
$ws.Range("E38").Value = "  -10.80%  "
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.30"
$ws.Range("E39").Value = "  -13.61%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0858"
$ws.Range("E40").Value = "  -14.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("E41").Value = "  -13.41%  "
$ws.Range("D42").Value = "2.902.82"
$ws.Range("E42").Value = "  -13.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  -24.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.265"
$ws.Range("E44").Value = "  -14.75%  "
$ws.Range("E45").Value = "  -22.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("E46").Value = "  -16.51%  "
$ws.Range("E47").Value = "  -15.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.57"
$ws.Range("E48").Value = "  -17.54%  "
$ws.Range("E50").Value = "  -12.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.93"
$ws.Range("E51").Value = "  -11.37%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
